$d = $word.ActiveDocument

# Locate the paragraph ending in "ha un suo codice." (the "Il Commit Hash..." paragraph)
# and insert a brand-new paragraph right after it, before the trailing blank paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ha un suo codice.*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate anchor paragraph (Il Commit Hash...)"
}

$targetIndex = $target.Index

$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newParaText = 'Banches sono le diramazioni del nostro progetto nella nostra Repository. Il ramo principale, detto anche ramo Master, è quello lineare del nostro progetto (es. se il nostro progetto riguarda i pesci e continuo a committare cose sui pesci sto lavorando sempre nel ramo Master). Se voglio lavorare allo stesso documento del ramo Master ma con uno scopo/obbiettivo diverso, divido il ramo principale o creo un ramo alternativo in qualunque punto dei miei Commit. Per fare questo in GitHub, entro nella mia Repository e creo un nuovo Branch (dandogli un nome), poi posso eseguire le modifiche ai file che voglio e committare il tutto. Quando committo posso decidere se farlo nel Branch appena creato o in un nuovo Branch (ovviamente posso anche commentare quello che ho appena modificato). Ora ho due versioni del file in due rami diversi: uno è nel ramo Master dove non ho le modifiche fatte nell’altro ramo; l’altro file è nel ramo appena creato con le modifiche appena effettuate (es se ho il mio file sui pesci e voglio parlare di cani ma non voglio modificare il file principale, basta creare un nuovo Branch chiamato “Cani”, effettuare le modifiche che voglio, committare ed ho finito. Se vado a vedere il file nel ramo Master troverò che parla solo di pesci, invece se vado a vederlo nel Branch “Cani” troverò che parla di cani e di pesci.). Se vado nella sezione Insights di GitHub posso controllare la mia diramazione (Branch) cliccando sull’opzione Network. Se con la mia diramazione voglio tornare nel ramo Master perché il lavoro che ho fatto può essere aggiunto al progetto principale, devo fare un Marge. Fare un Marge vine preceduto da una richiesta di Pull cioè una Pull Request. Sempre in GitHub trovo la sezione, nel mio Repository, Pull Request. Se clicco lì mi dice che c’è un nuovo Branch e posso comparare le modifiche nel Branch con quelle nel ramo Master. Una volta comparate e visto che non ci sono errori fa una richiesta di Pull e mi dice che non ci sono conflitti e che posso fare un Merge. Facendo un Merge unisco il ramo che si era distaccato dal Master con il ramo Master.'

$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.InsertAfter($newParaText)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
